$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-30 17:16:40"
$zhcn.Range("H2").Value = "2016-08-30 17:16:35"
$zhcn.Range("K2").Value = "2016-08-30 17:16:53"
$dede.Range("H2").Value = "2016-08-30 17:16:40"
$dede.Range("K2").Value = "2016-08-30 17:17:00"
